$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (A1:J1) to title-cased / human readable labels.
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Description"
$ws.Range("D1").Value = "Category"
$ws.Range("E1").Value = "Subcategory"
$ws.Range("F1").Value = "Country of origin"
$ws.Range("G1").Value = "Gross mass"
$ws.Range("H1").Value = "Net mass"
$ws.Range("I1").Value = "Weight unit"
$ws.Range("J1").Value = "Customs territories"

# Move the active cell selection to B2.
$ws.Range("B2").Select()
